$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Clear the Results value in E2 (was "PASS") while keeping its style,
# marking the test for (re-)execution.
$ws.Range("E2").ClearContents()

# Reflect the cell the user left selected after making the edit.
[void]$ws.Range("D7").Select()
